$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("draftpicks")

# 1. Correct the player name typo: "Josh James" -> "Joshua James" (row 427)
$ws.Range("B427").Value = "Joshua James"

# 2. Append new bench-draft pick rows (483-523): columns A, C, D, E first
#    (these never introduce new shared strings, order is not significant)
$ws.Range("A483").Value = "dembums"
$ws.Range("C483").Value = 0
$ws.Range("D483").Value = "B"
$ws.Range("E483").Value = 43878
$ws.Range("A484").Value = "chicago"
$ws.Range("C484").Value = 0
$ws.Range("D484").Value = "B"
$ws.Range("E484").Value = 43878
$ws.Range("A485").Value = "balco"
$ws.Range("C485").Value = 0
$ws.Range("D485").Value = "B"
$ws.Range("E485").Value = 43878
$ws.Range("A486").Value = "ds9"
$ws.Range("C486").Value = 0
$ws.Range("D486").Value = "B"
$ws.Range("E486").Value = 43878
$ws.Range("A487").Value = "dsb"
$ws.Range("C487").Value = 0
$ws.Range("D487").Value = "B"
$ws.Range("E487").Value = 43878
$ws.Range("A488").Value = "deener"
$ws.Range("C488").Value = 0
$ws.Range("D488").Value = "B"
$ws.Range("E488").Value = 43878
$ws.Range("A489").Value = "marmaduke"
$ws.Range("C489").Value = 0
$ws.Range("D489").Value = "B"
$ws.Range("E489").Value = 43878
$ws.Range("A490").Value = "sturgeon"
$ws.Range("C490").Value = 0
$ws.Range("D490").Value = "B"
$ws.Range("E490").Value = 43878
$ws.Range("A491").Value = "rippe"
$ws.Range("C491").Value = 0
$ws.Range("D491").Value = "B"
$ws.Range("E491").Value = 43879
$ws.Range("A492").Value = "sturgeon"
$ws.Range("C492").Value = 0
$ws.Range("D492").Value = "B"
$ws.Range("E492").Value = 43879
$ws.Range("A493").Value = "rippe"
$ws.Range("C493").Value = 0
$ws.Range("D493").Value = "B"
$ws.Range("E493").Value = 43879
$ws.Range("A494").Value = "ds9"
$ws.Range("C494").Value = 0
$ws.Range("D494").Value = "B"
$ws.Range("E494").Value = 43879
$ws.Range("A495").Value = "deano"
$ws.Range("C495").Value = 0
$ws.Range("D495").Value = "B"
$ws.Range("E495").Value = 43879
$ws.Range("A496").Value = "bears"
$ws.Range("C496").Value = 0
$ws.Range("D496").Value = "B"
$ws.Range("E496").Value = 43879
$ws.Range("A497").Value = "ottawa"
$ws.Range("C497").Value = 0
$ws.Range("D497").Value = "B"
$ws.Range("E497").Value = 43879
$ws.Range("A498").Value = "pasadena"
$ws.Range("C498").Value = 0
$ws.Range("D498").Value = "B"
$ws.Range("E498").Value = 43879
$ws.Range("A499").Value = "drjames"
$ws.Range("C499").Value = 0
$ws.Range("D499").Value = "B"
$ws.Range("E499").Value = 43879
$ws.Range("A500").Value = "balco"
$ws.Range("C500").Value = 0
$ws.Range("D500").Value = "B"
$ws.Range("E500").Value = 43879
$ws.Range("A501").Value = "rippe"
$ws.Range("C501").Value = 0
$ws.Range("D501").Value = "B"
$ws.Range("E501").Value = 43879
$ws.Range("A502").Value = "rippe"
$ws.Range("C502").Value = 0
$ws.Range("D502").Value = "B"
$ws.Range("E502").Value = 43879
$ws.Range("A503").Value = "drjames"
$ws.Range("C503").Value = 0
$ws.Range("D503").Value = "B"
$ws.Range("E503").Value = 43879
$ws.Range("A504").Value = "chicago"
$ws.Range("C504").Value = 0
$ws.Range("D504").Value = "B"
$ws.Range("E504").Value = 43879
$ws.Range("A505").Value = "dembums"
$ws.Range("C505").Value = 0
$ws.Range("D505").Value = "B"
$ws.Range("E505").Value = 43879
$ws.Range("A506").Value = "dembums"
$ws.Range("C506").Value = 0
$ws.Range("D506").Value = "B"
$ws.Range("E506").Value = 43880
$ws.Range("A507").Value = "chicago"
$ws.Range("C507").Value = 0
$ws.Range("D507").Value = "B"
$ws.Range("E507").Value = 43880
$ws.Range("A508").Value = "ds9"
$ws.Range("C508").Value = 0
$ws.Range("D508").Value = "B"
$ws.Range("E508").Value = 43880
$ws.Range("A509").Value = "marmaduke"
$ws.Range("C509").Value = 0
$ws.Range("D509").Value = "B"
$ws.Range("E509").Value = 43880
$ws.Range("A510").Value = "rippe"
$ws.Range("C510").Value = 0
$ws.Range("D510").Value = "B"
$ws.Range("E510").Value = 43880
$ws.Range("A511").Value = "balco"
$ws.Range("C511").Value = 0
$ws.Range("D511").Value = "B"
$ws.Range("E511").Value = 43880
$ws.Range("A512").Value = "drjames"
$ws.Range("C512").Value = 0
$ws.Range("D512").Value = "B"
$ws.Range("E512").Value = 43880
$ws.Range("A513").Value = "rippe"
$ws.Range("C513").Value = 0
$ws.Range("D513").Value = "B"
$ws.Range("E513").Value = 43880
$ws.Range("A514").Value = "ottawa"
$ws.Range("C514").Value = 0
$ws.Range("D514").Value = "B"
$ws.Range("E514").Value = 43880
$ws.Range("A515").Value = "bears"
$ws.Range("C515").Value = 0
$ws.Range("D515").Value = "B"
$ws.Range("E515").Value = 43880
$ws.Range("A516").Value = "bellevegas"
$ws.Range("C516").Value = 0
$ws.Range("D516").Value = "B"
$ws.Range("E516").Value = 43880
$ws.Range("A517").Value = "deano"
$ws.Range("C517").Value = 0
$ws.Range("D517").Value = "B"
$ws.Range("E517").Value = 43880
$ws.Range("A518").Value = "rippe"
$ws.Range("C518").Value = 0
$ws.Range("D518").Value = "B"
$ws.Range("E518").Value = 43880
$ws.Range("A519").Value = "drjames"
$ws.Range("C519").Value = 0
$ws.Range("D519").Value = "B"
$ws.Range("E519").Value = 43880
$ws.Range("A520").Value = "pkdodgers"
$ws.Range("C520").Value = 0
$ws.Range("D520").Value = "B"
$ws.Range("E520").Value = 43880
$ws.Range("A521").Value = "marmaduke"
$ws.Range("C521").Value = 0
$ws.Range("D521").Value = "B"
$ws.Range("E521").Value = 43880
$ws.Range("A522").Value = "pasadena"
$ws.Range("C522").Value = 0
$ws.Range("D522").Value = "B"
$ws.Range("E522").Value = 43880
$ws.Range("A523").Value = "deener"
$ws.Range("C523").Value = 0
$ws.Range("D523").Value = "B"
$ws.Range("E523").Value = 43880

# 3. Fill in column B (player names) in the exact order the names were
#    originally typed, so new shared-string entries are appended in the
#    same sequence as the source workbook.
$ws.Range("B483").Value = "Christian Pache"
$ws.Range("B484").Value = "Oneil Cruz"
$ws.Range("B485").Value = "CJ Abrams"
$ws.Range("B486").Value = "Jurickson Profar"
$ws.Range("B487").Value = "Miguel Cabrera"
$ws.Range("B488").Value = "Rowdy Tellez"
$ws.Range("B489").Value = "Drew Waters"
$ws.Range("B490").Value = "Adam Haseley"
$ws.Range("B492").Value = "Nicky Lopez"
$ws.Range("B494").Value = "Emmanuel Clase"
$ws.Range("B496").Value = "Bobby Witt Jr."
$ws.Range("B497").Value = "Franklin Barreto"
$ws.Range("B498").Value = "Brendan Rodgers"
$ws.Range("B499").Value = "Zach Plesac"
$ws.Range("B500").Value = "Kyle Wright"
$ws.Range("B501").Value = "Blake Treinen"
$ws.Range("B502").Value = "JJ Bleday"
$ws.Range("B503").Value = "Sheldon Neuse"
$ws.Range("B504").Value = "Alex Gordon"
$ws.Range("B491").Value = "Nico Hoerner"
$ws.Range("B505").Value = "Adbert Alzolay"
$ws.Range("B522").Value = "PLACEHOLDER"
$ws.Range("B523").Value = "PLACEHOLDER"
$ws.Range("B506").Value = "Jose Alvarado"
$ws.Range("B507").Value = "Lewis Thorpe"
$ws.Range("B508").Value = "Jeter Downs"
$ws.Range("B509").Value = "Spencer Turnbull"
$ws.Range("B510").Value = "Asdrubal Cabrera"
$ws.Range("B511").Value = "Nolan Jones"
$ws.Range("B512").Value = "Austin Voth"
$ws.Range("B513").Value = "Brian Goodwin"
$ws.Range("B514").Value = "Mitch Moreland"
$ws.Range("B515").Value = "Isan Diaz"
$ws.Range("B516").Value = "Leury Garcia"
$ws.Range("B517").Value = "Nick Markakis"
$ws.Range("B518").Value = "Shaun Anderson"
$ws.Range("B519").Value = "Lane Thomas"
$ws.Range("B520").Value = "Josiah Gray"
$ws.Range("B521").Value = "Corbin Burnes"
$ws.Range("B493").Value = "Garrett Cooper"
$ws.Range("B495").Value = "Collin McHugh"
